$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 299 (pushes the existing row 299..355 data down to 300..356)
$ws.Rows.Item(299).Insert()

# Populate the new row 299 with a new data record (weekly price update for Acelga)
$ws.Cells.Item(299, 1).Value = 4
$ws.Cells.Item(299, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(299, 3).Value = "Los Lagos"
$ws.Cells.Item(299, 4).Value = 45209
$ws.Cells.Item(299, 5).Value = 10
$ws.Cells.Item(299, 6).Value = 100112009
$ws.Cells.Item(299, 7).Value = "Acelga"
$ws.Cells.Item(299, 8).Value = "Sin especificar"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 75
$ws.Cells.Item(299, 11).Value = 10000
$ws.Cells.Item(299, 12).Value = 10000
$ws.Cells.Item(299, 13).Value = 10000
$ws.Cells.Item(299, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(299, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(299, 16).Value = 833
$ws.Cells.Item(299, 17).Value = 12
$ws.Cells.Item(299, 18).Value = "Hortaliza"

# Give the new date cell (column D) the same number format/style as the rest of
# the date column (the other rows' D cells carry a date/time number format).
$ws.Cells.Item(299, 4).NumberFormat = $ws.Cells.Item(300, 4).NumberFormat()
